# Add Quantity/Price columns (F, G) to the PRODUCTS sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---------------------------------------------------
$ws.Range("F1").Value = "Quantity"
$ws.Range("G1").Value = "Price"

# Match the header styling used by the existing D1/E1 header cells
# (bold white font, colored fill, border) by copying formats only.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2-31): Quantity (F), Price (G) --------------------------
$qty = @(0,0,0,0,0,0,0,0,0,0,0,0,2,0,0,0,1,0,0,0,0,0,0,1,0,0,0,0,0,0)
$price = @(1999,3299,1799,999,1999,449,499.99,799.99,0,0,0,25.99,15.99,13.99,23.99,23.99,23.99,23.99,33.450000000000003,23.99,48.59,23.99,237,523,229.95,49.95,599,149,828,9.99)

for ($i = 0; $i -lt 30; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $qty[$i]
    $ws.Cells.Item($row, 7).Value = $price[$i]
}

# --- Selection, matching the post-edit cursor position ------------------
$ws.Range("G1").Select() | Out-Null
